$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F2").Value = "'0.887207"
$ws.Range("G2").Value = "'0.869436"
$ws.Range("H2").Value = "'0.905719"
$ws.Range("I2").Value = "'0.960396"
$ws.Range("J2").Value = "'0.97"
$ws.Range("K2").Value = "'0.95098"
$ws.Range("L2").Value = "'0.623188"
$ws.Range("M2").Value = "'0.641791"
$ws.Range("N2").Value = "'0.605634"
$ws.Range("O2").Value = "'0.90625"
$ws.Range("P2").Value = "'0.988636"
$ws.Range("Q2").Value = "'0.836538"
$ws.Range("R2").Value = "'0.939597"
$ws.Range("S2").Value = "'0.985915"
$ws.Range("T2").Value = "'0.897436"
$ws.Range("U2").Value = "'0.966361"
$ws.Range("V2").Value = "'0.946108"
$ws.Range("W2").Value = "'0.9875"
$ws.Range("X2").Value = "'0.825397"
$ws.Range("Y2").Value = "'0.710383"
$ws.Range("Z2").Value = "'0.984848"
$ws.Range("F2:Z2").ClearFormats()
$ws.Range("F4").Value = "'0.882923"
$ws.Range("G4").Value = "'0.853026"
$ws.Range("H4").Value = "'0.914992"
$ws.Range("I4").Value = "'0.935961"
$ws.Range("J4").Value = "'0.940594"
$ws.Range("K4").Value = "'0.931373"
$ws.Range("L4").Value = "'0.647399"
$ws.Range("M4").Value = "'0.54902"
$ws.Range("N4").Value = "'0.788732"
$ws.Range("O4").Value = "'0.88172"
$ws.Range("P4").Value = "'1.0"
$ws.Range("Q4").Value = "'0.788462"
$ws.Range("R4").Value = "'0.921053"
$ws.Range("S4").Value = "'0.945946"
$ws.Range("T4").Value = "'0.897436"
$ws.Range("U4").Value = "'0.98452"
$ws.Range("V4").Value = "'0.97546"
$ws.Range("W4").Value = "'0.99375"
$ws.Range("X4").Value = "'0.855263"
$ws.Range("Y4").Value = "'0.755814"
$ws.Range("Z4").Value = "'0.984848"
$ws.Range("F4:Z4").ClearFormats()
$ws.Range("F6").Value = "'0.911677"
$ws.Range("G6").Value = "'0.88389"
$ws.Range("H6").Value = "'0.941267"
$ws.Range("I6").Value = "'0.970874"
$ws.Range("J6").Value = "'0.961538"
$ws.Range("K6").Value = "'0.980392"
$ws.Range("L6").Value = "'0.724832"
$ws.Range("M6").Value = "'0.692308"
$ws.Range("N6").Value = "'0.760563"
$ws.Range("O6").Value = "'0.96"
$ws.Range("P6").Value = "'1.0"
$ws.Range("Q6").Value = "'0.923077"
$ws.Range("R6").Value = "'0.939597"
$ws.Range("S6").Value = "'0.985915"
$ws.Range("T6").Value = "'0.897436"
$ws.Range("U6").Value = "'0.969325"
$ws.Range("V6").Value = "'0.951807"
$ws.Range("W6").Value = "'0.9875"
$ws.Range("X6").Value = "'0.846906"
$ws.Range("Y6").Value = "'0.742857"
$ws.Range("Z6").Value = "'0.984848"
$ws.Range("F6:Z6").ClearFormats()
$ws.Range("F8").Value = "'0.876712"
$ws.Range("G8").Value = "'0.863568"
$ws.Range("H8").Value = "'0.890263"
$ws.Range("I8").Value = "'0.952381"
$ws.Range("J8").Value = "'0.925926"
$ws.Range("K8").Value = "'0.980392"
$ws.Range("L8").Value = "'0.680556"
$ws.Range("M8").Value = "'0.671233"
$ws.Range("N8").Value = "'0.690141"
$ws.Range("O8").Value = "'0.88172"
$ws.Range("P8").Value = "'1.0"
$ws.Range("Q8").Value = "'0.788462"
$ws.Range("R8").Value = "'0.938776"
$ws.Range("S8").Value = "'1.0"
$ws.Range("T8").Value = "'0.884615"
$ws.Range("U8").Value = "'0.935065"
$ws.Range("V8").Value = "'0.972973"
$ws.Range("W8").Value = "'0.9"
$ws.Range("X8").Value = "'0.821317"
$ws.Range("Y8").Value = "'0.700535"
$ws.Range("Z8").Value = "'0.992424"
$ws.Range("F8:Z8").ClearFormats()
$ws.Range("F12").Value = "'0.878856"
$ws.Range("G12").Value = "'0.856305"
$ws.Range("H12").Value = "'0.902628"
$ws.Range("I12").Value = "'0.946341"
$ws.Range("J12").Value = "'0.941748"
$ws.Range("K12").Value = "'0.95098"
$ws.Range("L12").Value = "'0.580247"
$ws.Range("M12").Value = "'0.516484"
$ws.Range("N12").Value = "'0.661972"
$ws.Range("O12").Value = "'0.861702"
$ws.Range("P12").Value = "'0.964286"
$ws.Range("Q12").Value = "'0.778846"
$ws.Range("R12").Value = "'0.939597"
$ws.Range("S12").Value = "'0.985915"
$ws.Range("T12").Value = "'0.897436"
$ws.Range("U12").Value = "'0.972308"
$ws.Range("V12").Value = "'0.957576"
$ws.Range("W12").Value = "'0.9875"
$ws.Range("X12").Value = "'0.863787"
$ws.Range("Y12").Value = "'0.769231"
$ws.Range("Z12").Value = "'0.984848"
$ws.Range("F12:Z12").ClearFormats()
$ws.Range("F16").Value = "'0.900672"
$ws.Range("G16").Value = "'0.871387"
$ws.Range("H16").Value = "'0.931994"
$ws.Range("I16").Value = "'0.960396"
$ws.Range("J16").Value = "'0.97"
$ws.Range("K16").Value = "'0.95098"
$ws.Range("L16").Value = "'0.690058"
$ws.Range("M16").Value = "'0.59"
$ws.Range("N16").Value = "'0.830986"
$ws.Range("O16").Value = "'0.89899"
$ws.Range("P16").Value = "'0.946809"
$ws.Range("Q16").Value = "'0.855769"
$ws.Range("R16").Value = "'0.932432"
$ws.Range("S16").Value = "'0.985714"
$ws.Range("T16").Value = "'0.884615"
$ws.Range("U16").Value = "'0.9875"
$ws.Range("V16").Value = "'0.9875"
$ws.Range("W16").Value = "'0.9875"
$ws.Range("X16").Value = "'0.870432"
$ws.Range("Y16").Value = "'0.775148"
$ws.Range("Z16").Value = "'0.992424"
$ws.Range("F16:Z16").ClearFormats()
$ws.Range("F19").Value = "'0.884323"
$ws.Range("G19").Value = "'0.871064"
$ws.Range("H19").Value = "'0.897991"
$ws.Range("I19").Value = "'0.945274"
$ws.Range("J19").Value = "'0.959596"
$ws.Range("K19").Value = "'0.931373"
$ws.Range("L19").Value = "'0.661972"
$ws.Range("M19").Value = "'0.661972"
$ws.Range("N19").Value = "'0.661972"
$ws.Range("O19").Value = "'0.860104"
$ws.Range("P19").Value = "'0.932584"
$ws.Range("Q19").Value = "'0.798077"
$ws.Range("R19").Value = "'0.939597"
$ws.Range("S19").Value = "'0.985915"
$ws.Range("T19").Value = "'0.897436"
$ws.Range("U19").Value = "'0.971963"
$ws.Range("V19").Value = "'0.968944"
$ws.Range("W19").Value = "'0.975"
$ws.Range("X19").Value = "'0.83871"
$ws.Range("Y19").Value = "'0.730337"
$ws.Range("Z19").Value = "'0.984848"
$ws.Range("F19:Z19").ClearFormats()
$ws.Range("F22").Value = "'0.908277"
$ws.Range("G22").Value = "'0.877522"
$ws.Range("H22").Value = "'0.941267"
$ws.Range("I22").Value = "'0.959596"
$ws.Range("J22").Value = "'0.989583"
$ws.Range("K22").Value = "'0.931373"
$ws.Range("L22").Value = "'0.834356"
$ws.Range("M22").Value = "'0.73913"
$ws.Range("N22").Value = "'0.957746"
$ws.Range("O22").Value = "'0.90625"
$ws.Range("P22").Value = "'0.988636"
$ws.Range("Q22").Value = "'0.836538"
$ws.Range("R22").Value = "'0.939597"
$ws.Range("S22").Value = "'0.985915"
$ws.Range("T22").Value = "'0.897436"
$ws.Range("U22").Value = "'0.963415"
$ws.Range("V22").Value = "'0.940476"
$ws.Range("W22").Value = "'0.9875"
$ws.Range("X22").Value = "'0.83871"
$ws.Range("Y22").Value = "'0.730337"
$ws.Range("Z22").Value = "'0.984848"
$ws.Range("F22:Z22").ClearFormats()
$ws.Range("F25").Value = "'0.890226"
$ws.Range("G25").Value = "'0.866764"
$ws.Range("H25").Value = "'0.914992"
$ws.Range("I25").Value = "'0.971154"
$ws.Range("J25").Value = "'0.95283"
$ws.Range("K25").Value = "'0.990196"
$ws.Range("L25").Value = "'0.643836"
$ws.Range("M25").Value = "'0.626667"
$ws.Range("N25").Value = "'0.661972"
$ws.Range("O25").Value = "'0.895833"
$ws.Range("P25").Value = "'0.977273"
$ws.Range("Q25").Value = "'0.826923"
$ws.Range("R25").Value = "'0.939597"
$ws.Range("S25").Value = "'0.985915"
$ws.Range("T25").Value = "'0.897436"
$ws.Range("U25").Value = "'0.957576"
$ws.Range("V25").Value = "'0.929412"
$ws.Range("W25").Value = "'0.9875"
$ws.Range("X25").Value = "'0.852459"
$ws.Range("Y25").Value = "'0.751445"
$ws.Range("Z25").Value = "'0.984848"
$ws.Range("F25:Z25").ClearFormats()
$ws.Range("F26").Value = "'0.907186"
$ws.Range("G26").Value = "'0.879536"
$ws.Range("H26").Value = "'0.936631"
$ws.Range("I26").Value = "'0.957346"
$ws.Range("J26").Value = "'0.926606"
$ws.Range("K26").Value = "'0.990196"
$ws.Range("L26").Value = "'0.766234"
$ws.Range("M26").Value = "'0.710843"
$ws.Range("N26").Value = "'0.830986"
$ws.Range("O26").Value = "'0.910995"
$ws.Range("P26").Value = "'1.0"
$ws.Range("Q26").Value = "'0.836538"
$ws.Range("R26").Value = "'0.927152"
$ws.Range("S26").Value = "'0.958904"
$ws.Range("T26").Value = "'0.897436"
$ws.Range("U26").Value = "'0.957576"
$ws.Range("V26").Value = "'0.929412"
$ws.Range("W26").Value = "'0.9875"
$ws.Range("X26").Value = "'0.876254"
$ws.Range("Y26").Value = "'0.784431"
$ws.Range("Z26").Value = "'0.992424"
$ws.Range("F26:Z26").ClearFormats()
$ws.Range("F29").Value = "'0.882615"
$ws.Range("G29").Value = "'0.849785"
$ws.Range("H29").Value = "'0.918083"
$ws.Range("I29").Value = "'0.909091"
$ws.Range("J29").Value = "'0.88785"
$ws.Range("K29").Value = "'0.931373"
$ws.Range("L29").Value = "'0.723926"
$ws.Range("M29").Value = "'0.641304"
$ws.Range("N29").Value = "'0.830986"
$ws.Range("O29").Value = "'0.888889"
$ws.Range("P29").Value = "'0.988235"
$ws.Range("Q29").Value = "'0.807692"
$ws.Range("R29").Value = "'0.915033"
$ws.Range("S29").Value = "'0.933333"
$ws.Range("T29").Value = "'0.897436"
$ws.Range("U29").Value = "'0.954128"
$ws.Range("V29").Value = "'0.934132"
$ws.Range("W29").Value = "'0.975"
$ws.Range("X29").Value = "'0.852459"
$ws.Range("Y29").Value = "'0.751445"
$ws.Range("Z29").Value = "'0.984848"
$ws.Range("F29:Z29").ClearFormats()
$ws.Range("F34").Value = "'0.886038"
$ws.Range("G34").Value = "'0.865782"
$ws.Range("H34").Value = "'0.907264"
$ws.Range("I34").Value = "'0.961905"
$ws.Range("J34").Value = "'0.935185"
$ws.Range("K34").Value = "'0.990196"
$ws.Range("L34").Value = "'0.692308"
$ws.Range("M34").Value = "'0.635294"
$ws.Range("N34").Value = "'0.760563"
$ws.Range("O34").Value = "'0.887701"
$ws.Range("P34").Value = "'1.0"
$ws.Range("Q34").Value = "'0.798077"
$ws.Range("R34").Value = "'0.945946"
$ws.Range("S34").Value = "'1.0"
$ws.Range("T34").Value = "'0.897436"
$ws.Range("U34").Value = "'0.968944"
$ws.Range("V34").Value = "'0.962963"
$ws.Range("W34").Value = "'0.975"
$ws.Range("X34").Value = "'0.863787"
$ws.Range("Y34").Value = "'0.769231"
$ws.Range("Z34").Value = "'0.984848"
$ws.Range("F34:Z34").ClearFormats()
$ws.Range("F36").Value = "'0.898876"
$ws.Range("G36").Value = "'0.872093"
$ws.Range("H36").Value = "'0.927357"
$ws.Range("I36").Value = "'0.936585"
$ws.Range("J36").Value = "'0.932039"
$ws.Range("K36").Value = "'0.941176"
$ws.Range("L36").Value = "'0.652778"
$ws.Range("M36").Value = "'0.643836"
$ws.Range("N36").Value = "'0.661972"
$ws.Range("O36").Value = "'0.970297"
$ws.Range("P36").Value = "'1.0"
$ws.Range("Q36").Value = "'0.942308"
$ws.Range("R36").Value = "'0.915033"
$ws.Range("S36").Value = "'0.933333"
$ws.Range("T36").Value = "'0.897436"
$ws.Range("U36").Value = "'0.971963"
$ws.Range("V36").Value = "'0.968944"
$ws.Range("W36").Value = "'0.975"
$ws.Range("X36").Value = "'0.83871"
$ws.Range("Y36").Value = "'0.730337"
$ws.Range("Z36").Value = "'0.984848"
$ws.Range("F36:Z36").ClearFormats()
$ws.Range("F41").Value = "'0.882308"
$ws.Range("G41").Value = "'0.867164"
$ws.Range("H41").Value = "'0.897991"
$ws.Range("I41").Value = "'0.936585"
$ws.Range("J41").Value = "'0.932039"
$ws.Range("K41").Value = "'0.941176"
$ws.Range("L41").Value = "'0.619718"
$ws.Range("M41").Value = "'0.619718"
$ws.Range("N41").Value = "'0.619718"
$ws.Range("O41").Value = "'0.887701"
$ws.Range("P41").Value = "'1.0"
$ws.Range("Q41").Value = "'0.798077"
$ws.Range("R41").Value = "'0.927152"
$ws.Range("S41").Value = "'0.958904"
$ws.Range("T41").Value = "'0.897436"
$ws.Range("U41").Value = "'0.957576"
$ws.Range("V41").Value = "'0.929412"
$ws.Range("W41").Value = "'0.9875"
$ws.Range("X41").Value = "'0.860927"
$ws.Range("Y41").Value = "'0.764706"
$ws.Range("Z41").Value = "'0.984848"
$ws.Range("F41:Z41").ClearFormats()
$ws.Range("F42").Value = "'0.869369"
$ws.Range("G42").Value = "'0.845255"
$ws.Range("H42").Value = "'0.8949"
$ws.Range("I42").Value = "'0.95"
$ws.Range("J42").Value = "'0.969388"
$ws.Range("K42").Value = "'0.931373"
$ws.Range("L42").Value = "'0.630872"
$ws.Range("M42").Value = "'0.602564"
$ws.Range("N42").Value = "'0.661972"
$ws.Range("O42").Value = "'0.824859"
$ws.Range("P42").Value = "'1.0"
$ws.Range("Q42").Value = "'0.701923"
$ws.Range("R42").Value = "'0.933333"
$ws.Range("S42").Value = "'0.972222"
$ws.Range("T42").Value = "'0.897436"
$ws.Range("U42").Value = "'0.954128"
$ws.Range("V42").Value = "'0.934132"
$ws.Range("W42").Value = "'0.975"
$ws.Range("X42").Value = "'0.8125"
$ws.Range("Y42").Value = "'0.691489"
$ws.Range("Z42").Value = "'0.984848"
$ws.Range("F42:Z42").ClearFormats()
$ws.Range("F48").Value = "'0.891074"
$ws.Range("G48").Value = "'0.872593"
$ws.Range("H48").Value = "'0.910355"
$ws.Range("I48").Value = "'0.950495"
$ws.Range("J48").Value = "'0.96"
$ws.Range("K48").Value = "'0.941176"
$ws.Range("L48").Value = "'0.642857"
$ws.Range("M48").Value = "'0.652174"
$ws.Range("N48").Value = "'0.633803"
$ws.Range("O48").Value = "'0.896907"
$ws.Range("P48").Value = "'0.966667"
$ws.Range("Q48").Value = "'0.836538"
$ws.Range("R48").Value = "'0.935897"
$ws.Range("S48").Value = "'0.935897"
$ws.Range("T48").Value = "'0.935897"
$ws.Range("U48").Value = "'0.960725"
$ws.Range("V48").Value = "'0.929825"
$ws.Range("W48").Value = "'0.99375"
$ws.Range("X48").Value = "'0.860927"
$ws.Range("Y48").Value = "'0.764706"
$ws.Range("Z48").Value = "'0.984848"
$ws.Range("F48:Z48").ClearFormats()
$ws.Range("F53").Value = "'0.864382"
$ws.Range("G53").Value = "'0.834532"
$ws.Range("H53").Value = "'0.896445"
$ws.Range("I53").Value = "'0.975369"
$ws.Range("J53").Value = "'0.980198"
$ws.Range("K53").Value = "'0.970588"
$ws.Range("L53").Value = "'0.509554"
$ws.Range("M53").Value = "'0.465116"
$ws.Range("N53").Value = "'0.56338"
$ws.Range("O53").Value = "'0.857143"
$ws.Range("P53").Value = "'0.913043"
$ws.Range("Q53").Value = "'0.807692"
$ws.Range("R53").Value = "'0.945946"
$ws.Range("S53").Value = "'1.0"
$ws.Range("T53").Value = "'0.897436"
$ws.Range("U53").Value = "'0.945455"
$ws.Range("V53").Value = "'0.917647"
$ws.Range("W53").Value = "'0.975"
$ws.Range("X53").Value = "'0.850649"
$ws.Range("Y53").Value = "'0.744318"
$ws.Range("Z53").Value = "'0.992424"
$ws.Range("F53:Z53").ClearFormats()
$ws.Range("F55").Value = "'0.861727"
$ws.Range("G55").Value = "'0.851964"
$ws.Range("H55").Value = "'0.871716"
$ws.Range("I55").Value = "'0.941176"
$ws.Range("J55").Value = "'0.941176"
$ws.Range("K55").Value = "'0.941176"
$ws.Range("L55").Value = "'0.362069"
$ws.Range("M55").Value = "'0.466667"
$ws.Range("N55").Value = "'0.295775"
$ws.Range("O55").Value = "'0.913706"
$ws.Range("P55").Value = "'0.967742"
$ws.Range("Q55").Value = "'0.865385"
$ws.Range("R55").Value = "'0.939597"
$ws.Range("S55").Value = "'0.985915"
$ws.Range("T55").Value = "'0.897436"
$ws.Range("U55").Value = "'0.925816"
$ws.Range("V55").Value = "'0.881356"
$ws.Range("W55").Value = "'0.975"
$ws.Range("X55").Value = "'0.859016"
$ws.Range("Y55").Value = "'0.757225"
$ws.Range("Z55").Value = "'0.992424"
$ws.Range("F55:Z55").ClearFormats()
